{"js": "// Update the title date line and the 20x5 arithmetic practice table.\nconst body = context.document.body;\n\n// --- 1. Update the date/weekday title paragraph (first paragraph in body) ---\nconst paras = body.paragraphs;\nparas.load(\"items\");\nawait context.sync();\n\nconst titlePara = paras.items[0];\ntitlePara.load(\"text\");\nawait context.sync();\n\nif (titlePara.text === \"2024-08-04 Sunday\") {\n  titlePara.insertText(\"2024-08-05 Monday\", \"Replace\");\n}\n\n// --- 2. Update the table of arithmetic expressions ---\nconst tables = body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\nconst newValues = [\n  [\"86-24=\", \"34+32=\", \"6+72=\", \"22+7=\", \"5+20=\"],\n  [\"86-51=\", \"79-74=\", \"94-13=\", \"1+75=\", \"54-53=\"],\n  [\"86-17=\", \"59+29=\", \"27+56=\", \"94-68=\", \"42+43=\"],\n  [\"79+9=\", \"0+88=\", \"38+20=\", \"72-56=\", \"90-41=\"],\n  [\"39+9=\", \"67-50=\", \"0+99=\", \"19+34=\", \"59-44=\"],\n  [\"33-14=\", \"97-53=\", \"16+12=\", \"72-67=\", \"18+78=\"],\n  [\"35+12=\", \"74-74=\", \"49-40=\", \"85-61=\", \"45+42=\"],\n  [\"44-25=\", \"94-76=\", \"75+10=\", \"34+50=\", \"23+1=\"],\n  [\"56-37=\", \"15+62=\", \"59-7=\", \"68+30=\", \"41+41=\"],\n  [\"45+40=\", \"21+7=\", \"77-53=\", \"62-6=\", \"30+39=\"],\n  [\"65-60=\", \"34+6=\", \"48-29=\", \"15-10=\", \"41+26=\"],\n  [\"32-5=\", \"26+17=\", \"15+18=\", \"71-22=\", \"25-6=\"],\n  [\"2+18=\", \"9+41=\", \"84-72=\", \"62-43=\", \"38-22=\"],\n  [\"8+75=\", \"51-2=\", \"0+78=\", \"52+24=\", \"60-41=\"],\n  [\"47+10=\", \"62+18=\", \"17+19=\", \"67-61=\", \"28-15=\"],\n  [\"56-12=\", \"82-71=\", \"48-35=\", \"51+22=\", \"26-18=\"],\n  [\"80-37=\", \"93-76=\", \"69-39=\", \"16+49=\", \"45+40=\"],\n  [\"30-21=\", \"36+62=\", \"29+65=\", \"47-15=\", \"43-30=\"],\n  [\"82-33=\", \"33-27=\", \"78-28=\", \"71-33=\", \"47+12=\"],\n  [\"56+13=\", \"73-50=\", \"75-1=\", \"50-48=\", \"52-30=\"]\n];\n\ntable.values = newValues;\nawait context.sync();\n", "ps1": "# Update the title date line and the 20x5 arithmetic practice table.\n$d = $word.ActiveDocument\n\n# --- 1. Update the date/weekday title paragraph (first paragraph in doc) ---\n$titlePara = $d.Paragraphs(1)\nif ($titlePara.Range.Text.TrimEnd(\"`r\") -eq \"2024-08-04 Sunday\") {\n    $titlePara.Range.Text = \"2024-08-05 Monday\"\n}\n\n# --- 2. Update the table of arithmetic expressions ---\n$newValues = @(\n    @(\"86-24=\", \"34+32=\", \"6+72=\", \"22+7=\", \"5+20=\"),\n    @(\"86-51=\", \"79-74=\", \"94-13=\", \"1+75=\", \"54-53=\"),\n    @(\"86-17=\", \"59+29=\", \"27+56=\", \"94-68=\", \"42+43=\"),\n    @(\"79+9=\", \"0+88=\", \"38+20=\", \"72-56=\", \"90-41=\"),\n    @(\"39+9=\", \"67-50=\", \"0+99=\", \"19+34=\", \"59-44=\"),\n    @(\"33-14=\", \"97-53=\", \"16+12=\", \"72-67=\", \"18+78=\"),\n    @(\"35+12=\", \"74-74=\", \"49-40=\", \"85-61=\", \"45+42=\"),\n    @(\"44-25=\", \"94-76=\", \"75+10=\", \"34+50=\", \"23+1=\"),\n    @(\"56-37=\", \"15+62=\", \"59-7=\", \"68+30=\", \"41+41=\"),\n    @(\"45+40=\", \"21+7=\", \"77-53=\", \"62-6=\", \"30+39=\"),\n    @(\"65-60=\", \"34+6=\", \"48-29=\", \"15-10=\", \"41+26=\"),\n    @(\"32-5=\", \"26+17=\", \"15+18=\", \"71-22=\", \"25-6=\"),\n    @(\"2+18=\", \"9+41=\", \"84-72=\", \"62-43=\", \"38-22=\"),\n    @(\"8+75=\", \"51-2=\", \"0+78=\", \"52+24=\", \"60-41=\"),\n    @(\"47+10=\", \"62+18=\", \"17+19=\", \"67-61=\", \"28-15=\"),\n    @(\"56-12=\", \"82-71=\", \"48-35=\", \"51+22=\", \"26-18=\"),\n    @(\"80-37=\", \"93-76=\", \"69-39=\", \"16+49=\", \"45+40=\"),\n    @(\"30-21=\", \"36+62=\", \"29+65=\", \"47-15=\", \"43-30=\"),\n    @(\"82-33=\", \"33-27=\", \"78-28=\", \"71-33=\", \"47+12=\"),\n    @(\"56+13=\", \"73-50=\", \"75-1=\", \"50-48=\", \"52-30=\")\n)\n\n$t = $d.Tables(1)\nfor ($r = 1; $r -le 20; $r++) {\n    for ($c = 1; $c -le 5; $c++) {\n        $t.Cell($r, $c).Range.Text = $newValues[$r - 1][$c - 1]\n    }\n}\n"}
